$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# The newly inserted D:E columns inherit formatting from column C (left
# neighbor). Copy number formats/styles from F:G (the first two columns of
# the old data block, now shifted right) so D:E match the rest of the data
# columns in each row.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore explicit "best fit" column widths on the columns whose custom
# width metadata the column-insert left blank (the two brand-new columns,
# plus the columns that fell outside the original named col-ranges once
# everything shifted right).
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(7).ColumnWidth
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(12).ColumnWidth

# Populate the two new columns with the new quarter data (Dec-2018 / Sep-2018)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 22400
$ws.Range("E8").Value = 19600
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -400
$ws.Range("E15").Value = -300
$ws.Range("D17").Value = 4900
$ws.Range("E17").Value = 3800
$ws.Range("D18").Value = 17500
$ws.Range("E18").Value = 15800
$ws.Range("D20").Value = -11900
$ws.Range("E20").Value = -13100
$ws.Range("D21").Value = 6700
$ws.Range("E21").Value = 3700
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 5600
$ws.Range("E23").Value = 2700
$ws.Range("D24").Value = 900
$ws.Range("E24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 4600
$ws.Range("E26").Value = 2200
$ws.Range("D27").Value = 4600
$ws.Range("E27").Value = 2100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 11900
$ws.Range("E32").Value = 13100
$ws.Range("D33").Value = 4600
$ws.Range("E33").Value = 2100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 4600
$ws.Range("E35").Value = 2100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 24600
$ws.Range("E41").Value = 31100
$ws.Range("D42").Value = 15500
$ws.Range("E42").Value = 31000
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 25300
$ws.Range("E48").Value = 25500
$ws.Range("D49").Value = 70100
$ws.Range("E49").Value = 70500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 4700
$ws.Range("E52").Value = 4400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2078000
$ws.Range("E54").Value = 2044300
$ws.Range("D57").Value = 2300
$ws.Range("E57").Value = 2300
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 75100
$ws.Range("E61").Value = 45200
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1854800
$ws.Range("E66").Value = 1822400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 3400
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 39600
$ws.Range("E72").Value = 36200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 223200
$ws.Range("E76").Value = 218400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 4600
$ws.Range("E81").Value = 2100
$ws.Range("D83").Value = 1100
$ws.Range("E83").Value = 1000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 4700
$ws.Range("E89").Value = 1400
$ws.Range("D91").Value = -600
$ws.Range("E91").Value = -1100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -56900
$ws.Range("E94").Value = 34100
$ws.Range("D96").Value = -1300
$ws.Range("E96").Value = -1300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 30100
$ws.Range("E100").Value = -12800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -22000
$ws.Range("E102").Value = 22700

$ws.Range("A1").Select()
